# Fruta / hortaliza, semanal
#
# A new weekly price record is inserted as row 94 of the data table
# (pushing the existing rows 94-141 down to 95-142, dimension A1:R141 ->
# A1:R142).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 94, shifting rows 94:141
# down to 95:142 (keeps formatting such as the date style on column D).
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with the new weekly record.
$ws.Cells.Item(94, 1).Value  = 10
$ws.Cells.Item(94, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(94, 3).Value  = "La Araucanía"
$ws.Cells.Item(94, 4).Value  = 44879
$ws.Cells.Item(94, 5).Value  = 9
$ws.Cells.Item(94, 6).Value  = 100114002
$ws.Cells.Item(94, 7).Value  = "Camote"
$ws.Cells.Item(94, 8).Value  = "Sin especificar"
$ws.Cells.Item(94, 9).Value  = "Primera"
$ws.Cells.Item(94, 10).Value = 80
$ws.Cells.Item(94, 11).Value = 24000
$ws.Cells.Item(94, 12).Value = 24000
$ws.Cells.Item(94, 13).Value = 24000
$ws.Cells.Item(94, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(94, 15).Value = "Perú"
$ws.Cells.Item(94, 16).Value = 1200
$ws.Cells.Item(94, 17).Value = 20
$ws.Cells.Item(94, 18).Value = "Hortaliza"
